$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new data rows describing the card-game reader tables ---
$ws.Range("B7").Value = "cardGame.NormalCardReaderDB"
$ws.Range("C7").Value = "NormalCardDB"
$ws.Range("D7").Value = $true
$ws.Range("E7").Value = "cardGame/normalCardData.xlsx"

$ws.Range("B8").Value = "cardGame.EmotionCardReaderDB"
$ws.Range("C8").Value = "EmotionCardDB"
$ws.Range("D8").Value = $true
$ws.Range("E8").Value = "cardGame/emotionCardData.xlsx"

# --- Widen column B so the longer full_name values are readable ---
$ws.Columns.Item(2).ColumnWidth = 32.14

# --- Move the selection to where editing left off ---
$ws.Range("D11").Select()

# --- Bump the sheet's outline-level bookkeeping to track the two new
#     rows (mirrors the source file's existing row-count/outline-level
#     convention) without leaving a stray grouped row behind. ---
$ws.Rows.Item(100).OutlineLevel = 7
$ws.Rows.Item(100).Delete()
